$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ID" column (A) is being removed; it duplicated the "cliente_id"
# column that already existed at the far right (G) of the table, so the
# remaining columns (IBAN..clasificacion..cliente_id) shift one place to
# the left (B:G -> A:F).
#
# Do the shift through a non-overlapping staging area so every cell is
# moved (value + number format) exactly once, avoiding any read/write
# clobbering that an in-place overlapping copy would cause.
$stagingTopLeft = "A20"

$ws.Range("B1:G6").Copy($ws.Range($stagingTopLeft))
$ws.Range("A1:G6").Clear()
$ws.Range("A20:F25").Copy($ws.Range("A1"))
$ws.Range("A20:F25").Clear()

# The rectangular copy/paste leaves behind stub cells for positions that
# were blank in the source range; drop those so the sheet doesn't carry
# any empty <c> entries that weren't there before.
for ($r = 1; $r -le 6; $r++) {
  for ($c = 1; $c -le 6; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    if ($cell.Value2 -eq $null -or $cell.Value2 -eq "") {
      $cell.Clear()
    }
  }
}

# Leftover selection/formatting artifact: the user ended up with cell
# G12 selected after underlining it (same empty-but-underlined artifact
# already present at I6).
$ws.Range("G12").Font.Underline = $true
$ws.Range("G12").Select() | Out-Null
